$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 1.73
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.6
$ws.Range("Y2").Value = 19
$ws.Range("AW2").Value = 3.5
$ws.Range("BD5").Value = 151
$ws.Range("G6").Value = 2.4
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 3.2
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 2.2
$ws.Range("R6").Value = 1.65
$ws.Range("AA6").Value = 21
$ws.Range("AE6").Value = 15
$ws.Range("AG6").Value = 351
$ws.Range("AI6").Value = 15
$ws.Range("AJ6").Value = 12
$ws.Range("AK6").Value = 34
$ws.Range("G7").Value = 2.01
$ws.Range("J7").Value = 2.75
$ws.Range("L7").Value = 4.5
$ws.Range("X7").Value = 8.5
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 17
$ws.Range("AG7").Value = 451
$ws.Range("AH7").Value = 8.5
$ws.Range("AJ7").Value = 13
$ws.Range("AP7").Value = 23
$ws.Range("AZ7").Value = 81
$ws.Range("G8").Value = 3.5
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 2.2
$ws.Range("J8").Value = 4.33
$ws.Range("K8").Value = 1.91
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("W8").Value = 8
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 13
$ws.Range("AA8").Value = 34
$ws.Range("AB8").Value = 41
$ws.Range("AD8").Value = 6
$ws.Range("AE8").Value = 19
$ws.Range("AF8").Value = 67
$ws.Range("AH8").Value = 6
$ws.Range("AI8").Value = 9
$ws.Range("AJ8").Value = 10
$ws.Range("AK8").Value = 21
$ws.Range("AO8").Value = 21
$ws.Range("AP8").Value = 34
$ws.Range("AU8").Value = 9
$ws.Range("AV8").Value = 67
$ws.Range("AW8").Value = 4
$ws.Range("AX8").Value = 13
$ws.Range("I11").Value = 5
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("X11").Value = 7
$ws.Range("AD11").Value = 6.5
$ws.Range("AH11").Value = 12
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95
$ws.Range("BC13").Value = 151
$ws.Range("BD13").Value = 151
$ws.Range("G14").Value = 1.8
$ws.Range("I14").Value = 4.33
$ws.Range("J14").Value = 2.4
$ws.Range("L14").Value = 4.5
$ws.Range("O14").Value = 1.25
$ws.Range("P14").Value = 4
$ws.Range("Q14").Value = 1.83
$ws.Range("R14").Value = 2.03
$ws.Range("U14").Value = 1.73
$ws.Range("V14").Value = 2
$ws.Range("X14").Value = 9
$ws.Range("AB14").Value = 23
$ws.Range("AC14").Value = 12
$ws.Range("AG14").Value = 201
$ws.Range("AI14").Value = 21
$ws.Range("AJ14").Value = 13
$ws.Range("AK14").Value = 41
$ws.Range("AO14").Value = 9.5
$ws.Range("BA14").Value = 81
$ws.Range("Q15").Value = 2.2
$ws.Range("R15").Value = 1.67
$ws.Range("I16").Value = 4.1
$ws.Range("Q16").Value = 1.8
$ws.Range("R16").Value = 2
$ws.Range("W16").Value = 8.5
$ws.Range("Z16").Value = 17
$ws.Range("AD16").Value = 7
$ws.Range("AG16").Value = 151
$ws.Range("AJ16").Value = 13
$ws.Range("AL16").Value = 29
$ws.Range("AO16").Value = 10
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 8
$ws.Range("O17").Value = 1.36
$ws.Range("P17").Value = 3.2
$ws.Range("Q17").Value = 2.1
$ws.Range("R17").Value = 1.73
$ws.Range("G18").Value = 1.42
$ws.Range("H18").Value = 4.33
$ws.Range("J18").Value = 1.95
$ws.Range("K18").Value = 2.38
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 13
$ws.Range("O18").Value = 1.25
$ws.Range("P18").Value = 3.75
$ws.Range("Q18").Value = 1.83
$ws.Range("R18").Value = 2.03
$ws.Range("S18").Value = 1.36
$ws.Range("T18").Value = 3
$ws.Range("W18").Value = 7
$ws.Range("Y18").Value = 9
$ws.Range("AG18").Value = 401
$ws.Range("AH18").Value = 17
$ws.Range("AJ18").Value = 21
$ws.Range("AK18").Value = 81
$ws.Range("AQ18").Value = 21
$ws.Range("AT18").Value = 3
$ws.Range("AV18").Value = 67
$ws.Range("BB18").Value = 351
$ws.Range("G19").Value = 1.67
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 4.75
$ws.Range("O19").Value = 1.22
$ws.Range("P19").Value = 4
$ws.Range("U19").Value = 1.7
$ws.Range("V19").Value = 2.05
$ws.Range("X19").Value = 8.5
$ws.Range("AB19").Value = 23
$ws.Range("AG19").Value = 201
$ws.Range("AI19").Value = 26
$ws.Range("AY19").Value = 29
$ws.Range("G21").Value = 1.9
$ws.Range("H21").Value = 3.6
$ws.Range("I21").Value = 3.9
$ws.Range("O21").Value = 1.17
$ws.Range("P21").Value = 5
$ws.Range("Q21").Value = 1.6
$ws.Range("R21").Value = 2.3
$ws.Range("S21").Value = 1.3
$ws.Range("T21").Value = 3.4
$ws.Range("U21").Value = 1.57
$ws.Range("V21").Value = 2.25
$ws.Range("W21").Value = 9.5
$ws.Range("AD21").Value = 7.5
$ws.Range("AE21").Value = 12
$ws.Range("AF21").Value = 41
$ws.Range("AI21").Value = 21
$ws.Range("AL21").Value = 26
$ws.Range("AT21").Value = 3.4
$ws.Range("BC21").Value = 401
